# "fixing practice q code for array"
#
# Updates the "practiceQ" sheet:
#  - Fixes the buggy findMaxConsecutiveOnes() sample code (it used to bail
#    out of the loop on the first 0 because of a stray `return` inside the
#    `else` branch) and refreshes its expected result.
#  - Replaces the duplicated findMaxConsecutiveOnes code that had been
#    pasted into the "even number of digits" column with the real
#    findNumbers() solution, and refreshes its expected result.
#  - Renames the RunResult header to "Result" and fills in the previously
#    empty SubmitResult column (G1/G4).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("practiceQ")

# --- Row 1 headers: SubmitResult column now has a header ---
$ws.Range("G1").Value = "SubmitResult"

# --- Row 2: sample code cells ---
$findMaxCode = @'
def findMaxConsecutiveOnes(nums):
    max_count = 0
    current_count = 0
    for num in nums:
        if num == 1:
            current_count += 1
            max_count = max(max_count, current_count)
        else:
            current_count = 0  # reset when 0 is found
    return max_count
# Example usage:
print(findMaxConsecutiveOnes([1,1,0,1,1,1]))  # Output: 3
print(findMaxConsecutiveOnes([1,0,1,1,0,1]))  # Output: 2
'@
$ws.Range("C2").Value = $findMaxCode

$findNumbersCode = @'
def findNumbers(nums):
    count = 0
    for num in nums:
        if len(str(num)) % 2 == 0:
            count += 1
    return count
# Example usage:
print(findNumbers([12, 345, 2, 6, 7896]))   # Output: 2
print(findNumbers([555, 901, 482, 1771]))  # Output: 1
'@
$ws.Range("D2").Value = $findNumbersCode

# --- Row 4: expected results ---
$ws.Range("A4").Value = "Result"

$ws.Range("C4").Value = "3`n2"
$ws.Range("D4").Value = "2`n1"

# C4 used to carry a leftover "Arimo" font style; line it up with the
# "Play" font already used by its sibling result cell (B4).
$ws.Range("B4").Copy()
$ws.Range("C4").PasteSpecial(-4122)  # xlPasteFormats

# F4/G4 are new cells -- copy formats from matching cells elsewhere in
# the sheet (F2 is the same "hello" PythonCode column, B5 already has
# the Arial/top/wrap look used for SubmitResult cells).
$ws.Range("F2").Copy()
$ws.Range("F4").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F4").Value = "hello"

$ws.Range("B5").Copy()
$ws.Range("G4").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("G4").Value = "Submission Successful"

$excel.CutCopyMode = $false

# Selection ends on E2 after the edit.
$ws.Range("E2").Select()
